$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7
$ws.Range("G7").Value = 1.68
$ws.Range("P7").Value = 1.41
$ws.Range("Q7").Value = 2.62

# Row 8
$ws.Range("G8").Value = 2.92
$ws.Range("H8").Value = 3.05
$ws.Range("I8").Value = 2.37
$ws.Range("L8").Value = 1.32
$ws.Range("M8").Value = 2.82
$ws.Range("P8").Value = 1.44
$ws.Range("Q8").Value = 2.4
$ws.Range("R8").Value = 1.72
$ws.Range("S8").Value = 1.9
$ws.Range("T8").Value = 8.75
$ws.Range("U8").Value = 15
$ws.Range("W8").Value = 37
$ws.Range("Y8").Value = 35
$ws.Range("AA8").Value = 5.9
$ws.Range("AB8").Value = 13.5
$ws.Range("AC8").Value = 65
$ws.Range("AE8").Value = 7.6
$ws.Range("AH8").Value = 25
$ws.Range("AI8").Value = 20

# Row 9
$ws.Range("R9").Value = 2.05

# Row 14
$ws.Range("G14").Value = 1.62
$ws.Range("H14").Value = 3.7
$ws.Range("I14").Value = 4.7
$ws.Range("K14").Value = 7.5
$ws.Range("L14").Value = 1.29
$ws.Range("M14").Value = 3.3
$ws.Range("N14").Value = 1.85
$ws.Range("O14").Value = 1.85
$ws.Range("P14").Value = 1.39
$ws.Range("Q14").Value = 2.77
$ws.Range("R14").Value = 1.88
$ws.Range("S14").Value = 1.82
$ws.Range("T14").Value = 6.6
$ws.Range("U14").Value = 7.4
$ws.Range("W14").Value = 11.75
$ws.Range("X14").Value = 13.5
$ws.Range("Z14").Value = 7.5
$ws.Range("AA14").Value = 7.4
$ws.Range("AB14").Value = 17
$ws.Range("AE14").Value = 13
$ws.Range("AF14").Value = 27
$ws.Range("AG14").Value = 15.5
$ws.Range("AH14").Value = 80
$ws.Range("AI14").Value = 50
$ws.Range("AJ14").Value = 55
